$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "passwords" table renamed to "storagepass"
$ws.Range("A2").Value = "storagepass"

# 2. Insert two blank rows before the old "local" table block (row 6),
#    pushing it down to start at row 8.
$ws.Range("A6:A7").EntireRow.Insert()

# 3. Old "local" table renamed to "localconf" (now at row 8 after the insert)
$ws.Range("A8").Value = "localconf"

# 4. New field "value0" added ahead of the old value1..value5 list - the
#    existing rows are renumbered down by one (value1->value0, value2->value1, ...)
$ws.Range("B11").Value = "value0"
$ws.Range("B12").Value = "value1"
$ws.Range("B13").Value = "value2"
$ws.Range("B14").Value = "value3"
$ws.Range("B15").Value = "value4"
$ws.Range("C11").Value = "TEXT"
$ws.Range("C12").Value = "TEXT"
$ws.Range("C13").Value = "TEXT"
$ws.Range("C14").Value = "TEXT"
$ws.Range("C15").Value = "TEXT"

# 5. New comment for the password field
$ws.Range("D4").Value = "base64-encoded"

# 6. New note row under the storagepass table explaining the unique constraint
$ws.Range("A6").Value = "ALTER TABLE storagepass ADD CONSTRAINT onePasswordPerPool UNIQUE (vfs , pool)"
$ws.Range("A6").Interior.ThemeColor = 10

# 7. New note row under the localconf table explaining the unique constraint
$ws.Range("A17").Value = "ALTER TABLE localconf ADD CONSTRAINT oneValuesetPerKey UNIQUE (service , section , key)"
$ws.Range("A17").Interior.ThemeColor = 10

# Match the author's final selection before save
$ws.Range("A18").Select()
